$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in column B (stored as date serials, formatted as dates)
$ws.Range("B2").Value = 42811
$ws.Range("B3").Value = 42812
$ws.Range("B4").Value = 42813
$ws.Range("B5").Value = 42814
$ws.Range("B6").Value = 42815
$ws.Range("B7").Value = 42816
$ws.Range("B8").Value = 42817
$ws.Range("B9").Value = 42818
$ws.Range("B10").Value = 42819

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("B2:B10").Select()
